$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4333.3335
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 4333.3335
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 13000.0005
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -14748.0005
$ws.Range("H70").Value = 63163.418
$ws.Range("J70").Value = 83790.11
$ws.Range("L70").Value = 251370.33
$ws.Range("N70").Value = -251910.33
$ws.Range("H72").Value = 4333.3335
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 4333.3335
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 39000.0015
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -47736.0015
$ws.Range("H73").Value = 63163.418
$ws.Range("J73").Value = 83790.11
$ws.Range("L73").Value = 251370.33
$ws.Range("N73").Value = -253242.33
$ws.Range("H112").Value = 2373.6843
$ws.Range("I112").Value = 994.5
$ws.Range("J112").Value = 2535.9412
$ws.Range("K112").Value = 2983.5
$ws.Range("L112").Value = 7607.823600000001
$ws.Range("M112").Value = -1875.5
$ws.Range("N112").Value = -9823.8236
$ws.Range("H135").Value = 912.8421
$ws.Range("I135").Value = 637.94116
$ws.Range("K135").Value = 5741.47044
$ws.Range("M135").Value = -3206.47044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 907.5714
$ws.Range("I2").Value = 1238
$ws.Range("J2").Value = 467
$ws.Range("K2").Value = 1238
$ws.Range("L2").Value = 467
$ws.Range("M2").Value = -1125
$ws.Range("N2").Value = -693
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H45").Value = 1934.4
$ws.Range("I45").Value = 1844.25
$ws.Range("K45").Value = 1844.25
$ws.Range("M45").Value = -1467.25
$ws.Range("H63").Value = 6317.8823
$ws.Range("I63").Value = 5487.25
$ws.Range("K63").Value = 5487.25
$ws.Range("M63").Value = -4801.25
$ws.Range("H66").Value = 6317.8823
$ws.Range("I66").Value = 5487.25
$ws.Range("K66").Value = 27436.25
$ws.Range("M66").Value = -24004.25
$ws.Range("H102").Value = 448.45456
$ws.Range("J102").Value = 410
$ws.Range("L102").Value = 410
$ws.Range("N102").Value = -3654
$ws.Range("H110").Value = 3299.2
$ws.Range("I110").Value = 4602.857
$ws.Range("K110").Value = 4602.857
$ws.Range("M110").Value = -2557.857
$ws.Range("H116").Value = 907.5714
$ws.Range("I116").Value = 1238
$ws.Range("J116").Value = 467
$ws.Range("K116").Value = 1238
$ws.Range("L116").Value = 467
$ws.Range("M116").Value = 1056
$ws.Range("N116").Value = -5055
$ws.Range("H132").Value = 1686.0476
$ws.Range("I132").Value = 1600.3846
$ws.Range("K132").Value = 4801.1538
$ws.Range("M132").Value = -2271.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 907.5714
$ws.Range("I3").Value = 1238
$ws.Range("J3").Value = 467
$ws.Range("K3").Value = 1238
$ws.Range("L3").Value = 467
$ws.Range("M3").Value = -1124
$ws.Range("N3").Value = -695
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H20").Value = 858.5
$ws.Range("I20").Value = 898.2857
$ws.Range("J20").Value = 580
$ws.Range("K20").Value = 898.2857
$ws.Range("L20").Value = 580
$ws.Range("M20").Value = -651.2857
$ws.Range("N20").Value = -1074
$ws.Range("H64").Value = 733.125
$ws.Range("I64").Value = 562.4286
$ws.Range("K64").Value = 562.4286
$ws.Range("M64").Value = -337.4286
$ws.Range("H67").Value = 733.125
$ws.Range("I67").Value = 562.4286
$ws.Range("K67").Value = 562.4286
$ws.Range("M67").Value = 217.5714
$ws.Range("H86").Value = 1750
$ws.Range("I86").Value = 1750
$ws.Range("K86").Value = 1750
$ws.Range("M86").Value = -627
$ws.Range("H89").Value = 1750
$ws.Range("I89").Value = 1750
$ws.Range("K89").Value = 8750
$ws.Range("M89").Value = -3134
$ws.Range("H94").Value = 3039.1667
$ws.Range("I94").Value = 2445
$ws.Range("K94").Value = 2445
$ws.Range("M94").Value = -1994
$ws.Range("H99").Value = 4872.5557
$ws.Range("I99").Value = 4812.3335
$ws.Range("K99").Value = 4812.3335
$ws.Range("M99").Value = -3314.3335
$ws.Range("H105").Value = 4113.636
$ws.Range("I105").Value = 3073.158
$ws.Range("K105").Value = 3073.158
$ws.Range("M105").Value = -1326.158
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1718.375
$ws.Range("I16").Value = 949.75
$ws.Range("K16").Value = 949.75
$ws.Range("M16").Value = -662.75
$ws.Range("H31").Value = 5132.0713
$ws.Range("I31").Value = 1996.8
$ws.Range("J31").Value = 6873.8887
$ws.Range("K31").Value = 1996.8
$ws.Range("L31").Value = 6873.8887
$ws.Range("M31").Value = -1701.8
$ws.Range("N31").Value = -7463.8887
$ws.Range("H34").Value = 5132.0713
$ws.Range("I34").Value = 1996.8
$ws.Range("J34").Value = 6873.8887
$ws.Range("K34").Value = 1996.8
$ws.Range("L34").Value = 6873.8887
$ws.Range("M34").Value = -1794.8
$ws.Range("N34").Value = -7277.8887
$ws.Range("H52").Value = 96999.5
$ws.Range("J52").Value = 96999.5
$ws.Range("L52").Value = 96999.5
$ws.Range("N52").Value = -97587.5
$ws.Range("H62").Value = 82299.60000000001
$ws.Range("I62").Value = 2833
$ws.Range("J62").Value = 201499.5
$ws.Range("K62").Value = 2833
$ws.Range("L62").Value = 201499.5
$ws.Range("M62").Value = -2209
$ws.Range("N62").Value = -202747.5
$ws.Range("H65").Value = 82299.60000000001
$ws.Range("I65").Value = 2833
$ws.Range("J65").Value = 201499.5
$ws.Range("K65").Value = 14165
$ws.Range("L65").Value = 1007497.5
$ws.Range("M65").Value = -11045
$ws.Range("N65").Value = -1013737.5
$ws.Range("H107").Value = 726.5
$ws.Range("I107").Value = 521.73334
$ws.Range("J107").Value = 1165.2858
$ws.Range("K107").Value = 521.73334
$ws.Range("L107").Value = 1165.2858
$ws.Range("M107").Value = 1398.26666
$ws.Range("N107").Value = -5005.2858
$ws.Range("H113").Value = 1718.375
$ws.Range("I113").Value = 949.75
$ws.Range("K113").Value = 949.75
$ws.Range("M113").Value = 1220.25
$ws.Range("H132").Value = 2803.5789
$ws.Range("I132").Value = 1814.5454
$ws.Range("K132").Value = 5443.6362
$ws.Range("M132").Value = -2913.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 630.8857400000001
$ws.Range("J107").Value = 610.82355
$ws.Range("L107").Value = 1832.47065
$ws.Range("N107").Value = -5672.470649999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3997.25
$ws.Range("I113").Value = 2994.6667
$ws.Range("J113").Value = 4999.8335
$ws.Range("K113").Value = 2994.6667
$ws.Range("L113").Value = 4999.8335
$ws.Range("M113").Value = -824.6667000000002
$ws.Range("N113").Value = -9339.833500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 9854.714
$ws.Range("J122").Value = 10500
$ws.Range("L122").Value = 31500
$ws.Range("N122").Value = -36400
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2151.4546
$ws.Range("I132").Value = 1911.1666
$ws.Range("K132").Value = 5733.4998
$ws.Range("M132").Value = -3203.4998
$ws.Range("H136").Value = 2621.037
$ws.Range("I136").Value = 837.82355
$ws.Range("K136").Value = 2513.47065
$ws.Range("M136").Value = 36.52935000000025
